# Append new scrape results (2025-09-29 18:26 JST) to the "ランサーズ" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-29 18:26:32"

# --- 1. Make room for the two brand-new listings -----------------------
# Final layout: new job lands on row 5 (Snowflake Intelligence) and on
# row 9 (MySQL -> Google Sheets), pushing the previously scraped rows down.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(9).Insert()

# --- 2. Refresh the "取得日時" (fetched-at) timestamp on every row ------
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $newTimestamp
}

# --- 3. Row 5: brand-new listing ---------------------------------------
$ws.Cells.Item(5, 2).Value2 = "【急募】Snowflake IntelligenceでのAgent開発者を探しています!"
$ws.Cells.Item(5, 3).Value2 = "システム開発"
$ws.Cells.Item(5, 4).Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(5, 5).Value2 = "期限情報なし"
$ws.Cells.Item(5, 6).Value2 = "https://www.lancers.jp/work/detail/5403054"
$ws.Cells.Item(5, 7).Value2 = 68
$ws.Cells.Item(5, 8).Value2 = "◆開発"

# --- 4. Row 6: previously row 5 (Pawsitive) - price/score/skill updated
$ws.Cells.Item(6, 2).Value2 = "【ペットのアバター化】Pawsitiveプロトタイプ開発の依頼"
$ws.Cells.Item(6, 3).Value2 = "システム開発"
$ws.Cells.Item(6, 4).Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(6, 5).Value2 = "期限情報なし"
$ws.Cells.Item(6, 6).Value2 = "https://www.lancers.jp/work/detail/5399313"
$ws.Cells.Item(6, 7).Value2 = 68
$ws.Cells.Item(6, 8).Value2 = "◆開発"

# --- 5. Row 7: previously row 6 (kuchikomi site) - unchanged content ---
$ws.Cells.Item(7, 2).Value2 = "【急募】新しい口コミサイトの構築をお手伝いください!"
$ws.Cells.Item(7, 3).Value2 = "システム開発"
$ws.Cells.Item(7, 4).Value2 = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(7, 5).Value2 = "期限情報なし"
$ws.Cells.Item(7, 6).Value2 = "https://www.lancers.jp/work/detail/5402277"
$ws.Cells.Item(7, 7).Value2 = 38
$ws.Cells.Item(7, 8).Value2 = "◇サイト"

# --- 6. Row 8: previously row 7 (brand site alert) - unchanged content -
$ws.Cells.Item(8, 2).Value2 = "【急募】ブランドサイトの新商品更新アラート作成依頼"
$ws.Cells.Item(8, 3).Value2 = "システム開発"
$ws.Cells.Item(8, 4).Value2 = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(8, 5).Value2 = "期限情報なし"
$ws.Cells.Item(8, 6).Value2 = "https://www.lancers.jp/work/detail/5402794"
$ws.Cells.Item(8, 7).Value2 = 33
$ws.Cells.Item(8, 8).Value2 = "◇サイト"

# --- 7. Row 9: brand-new listing ----------------------------------------
$ws.Cells.Item(9, 2).Value2 = "MYSQLからGoogleスプレッドシートへデータ取り込み及びスプレッドシート改修"
$ws.Cells.Item(9, 3).Value2 = "システム開発"
$ws.Cells.Item(9, 4).Value2 = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(9, 5).Value2 = "期限情報なし"
$ws.Cells.Item(9, 6).Value2 = "https://www.lancers.jp/work/detail/5400606"
$ws.Cells.Item(9, 7).Value2 = 30
$ws.Cells.Item(9, 8).Value2 = "◇MySQL"

# --- 8. Row 10: previously row 8 (limited-disclosure job) --------------
$ws.Cells.Item(10, 2).Value2 = "限定公開 PR 限定公開の仕事"
$ws.Cells.Item(10, 3).Value2 = "システム開発"
$ws.Cells.Item(10, 4).Value2 = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(10, 5).Value2 = "期限情報なし"
$ws.Cells.Item(10, 6).Value2 = "https://www.lancers.jp/work/detail/5399347"
$ws.Cells.Item(10, 7).Value2 = 13

# --- 9. Row 11: previously row 9 (engineer interview) -------------------
$ws.Cells.Item(11, 2).Value2 = "エンジニア面談をお願い致します"
$ws.Cells.Item(11, 3).Value2 = "システム開発"
$ws.Cells.Item(11, 4).Value2 = "~ 5,000 円 / 固定"
$ws.Cells.Item(11, 5).Value2 = "期限情報なし"
$ws.Cells.Item(11, 6).Value2 = "https://www.lancers.jp/work/detail/5402603"
$ws.Cells.Item(11, 7).Value2 = 10

# --- 10. Rebuild the hyperlinks so the URL column (F) links correctly --
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 11; $r++) {
    $url = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $url)
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}

# --- 11. Widen column B to fit the new (longer) title text -------------
# (ColumnWidth uses Excel's character-width unit, which is offset by
#  ~5/6 of a character from the width persisted in the worksheet XML, so
#  back that padding out to land on an exact stored width of 46.)
$ws.Columns.Item(2).ColumnWidth = 45.16666666666666
